# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, matching the newly published site snapshot.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 283
$ws1.Range("F4").Value  = 1173
$ws1.Range("F5").Value  = 128
$ws1.Range("F6").Value  = 2790
$ws1.Range("F8").Value  = 710
$ws1.Range("F9").Value  = 122
$ws1.Range("F10").Value = 298
$ws1.Range("F12").Value = 714
$ws1.Range("F13").Value = 116
$ws1.Range("F15").Value = 1836
$ws1.Range("F16").Value = 312
$ws1.Range("F18").Value = 208

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 19
$ws2.Range("F7").Value  = 26
$ws2.Range("F10").Value = 59
$ws2.Range("F11").Value = 45
$ws2.Range("F18").Value = 37
$ws2.Range("F23").Value = 36

# 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6368
$ws3.Range("F4").Value = 2035
$ws3.Range("F5").Value = 275

# 全部类型 (All Types) - combined listing
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6368
$ws4.Range("F4").Value  = 2035
$ws4.Range("F5").Value  = 275
$ws4.Range("F11").Value = 283
$ws4.Range("F12").Value = 1173
$ws4.Range("F13").Value = 128
$ws4.Range("F14").Value = 19
$ws4.Range("F15").Value = 26
$ws4.Range("F17").Value = 2790
$ws4.Range("F20").Value = 59
$ws4.Range("F21").Value = 45
$ws4.Range("F23").Value = 710
$ws4.Range("F24").Value = 122
$ws4.Range("F25").Value = 298
$ws4.Range("F28").Value = 714
$ws4.Range("F29").Value = 116
$ws4.Range("F32").Value = 1836
$ws4.Range("F33").Value = 312
$ws4.Range("F37").Value = 208
$ws4.Range("F39").Value = 37
$ws4.Range("F44").Value = 36
